$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray border formatting that used to mark A5 as the last row
$ws.Range("A5").ClearFormats()

# Add the new "antigen_detection" rule row (row 6).
# Values are written in reverse column order (C, then B, then A) so that the
# workbook's shared-strings table ends up in the same order as the target file.
$ws.Range("C6").Value = "Yes, antigen detected;No, no antigen detection"
$ws.Range("B6").Value = "Was there a positive antigen test (e.g. PCR positive in synovia)?"
$ws.Range("A6").Value = "antigen_detection"

# Match the active selection left behind in the saved file
$ws.Range("A6").Select()
